$wb = $excel.ActiveWorkbook

# Add the new "OpenAccountTest" worksheet after the existing "AddCustomerTest" sheet.
$existing = $wb.Worksheets.Item("AddCustomerTest")
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $existing)
$newSheet.Name = "OpenAccountTest"

# Populate header row and data row.
$newSheet.Range("A1").Value = "customer"
$newSheet.Range("B1").Value = "currency"
$newSheet.Range("A2").Value = "Katya Smith"
$newSheet.Range("B2").Value = "Dollar"

# Selection on the new sheet.
[void]$newSheet.Range("E7").Select()
